# The commit swaps the contents of ppt/theme/theme1.xml and ppt/theme/theme2.xml:
#   theme1.xml goes from the standard "Office Theme" palette to the "Integral" palette
#   theme2.xml goes from the "Integral" palette to the standard "Office Theme" palette
# (file names/relationships are untouched - only the <a:clrScheme> colours inside the
# two theme parts are exchanged). theme2.xml is the presentation's live/primary theme
# part (referenced by presentation.xml and slideMaster1.xml), so we drive the change
# through the live ThemeColorScheme object, pushing it from "Integral" to the classic
# Office 12-colour palette that used to live in theme1.xml.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Office Theme palette (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink)
$officeHex = @("000000", "FFFFFF", "44546A", "E7E6E6", "5B9BD5", "ED7D31", "A5A5A5", "FFC000", "4472C4", "70AD47", "0563C1", "954F72")

for ($i = 1; $i -le 12; $i++) {
    $hex = $officeHex[$i - 1]
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    $tcs.Item($i).RGB = $r + ($g * 256) + ($b * 65536)
}
